$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the date/station cells for rows 60 and 61 (leave C/D/E empty cells as-is)
$ws.Range("A60:B61").Clear()

# Update the active selection on the sheet
$ws.Range("C65").Select()
